# Bill of Materials-Mfd3 update:
# - Replace IDC connector row with an LED 0603 row (qty 34)
# - Remove the USB-B (P5) connector row entirely
# - Replace the TS-H003 tactile switch line with a smaller 5.2mm tact switch
# - Adjust the print scale

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Step 1: remove the USB-B / P5 connector row (old row 10) ---------------
$ws.Rows.Item(10).Delete()

# --- Step 2: turn row 6 (previously the IDC 2X4 2.54 connector) into the LED row
$ws.Range("A6").Value = "LED 0603"
$ws.Range("B6").Value = "LED 0603"
$ws.Range("C6").Value = "LED1, LED2, LED3, LED4, LED5, LED6, LED7, LED8, LED9, LED10, LED11, LED12, LED13, LED14, LED15, LED16, LED17, LED18, LED19, LED20, LED21, LED22, LED23, LED24, LED25, LED26, LED27, LED28, LED29, LED30, LED31, LED32, LED33, LED34"
$ws.Range("D6").Value = "LED 0603"
$ws.Range("E6").Value = 34
$ws.Range("F6").Value = "C72043"

# Re-apply the original "text" cell format (lost when the .Value assignment
# rewrote the cell), using an untouched sibling cell as the format source.
$ws.Range("A9").Copy()
$ws.Range("A6:D6").PasteSpecial($xlPasteFormats)
$ws.Range("A9").Copy()
$ws.Range("F6").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(6).RowHeight = 180

# --- Step 3: update the tactile-switch row (now row 15 after the deletion) --
$ws.Range("A15").Value = "TACT SWITCH 5.2mm"
$ws.Range("B15").Value = "Push button 5.2 mm"
$ws.Range("D15").Value = "TACT SWITCH 5.2mm"
$ws.Range("F15").Value = "C412369"

$ws.Range("C9").Copy()
$ws.Range("A15:B15").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Copy()
$ws.Range("F15").PasteSpecial($xlPasteFormats)

# --- Step 4: adjust the print scale -----------------------------------------
$ws.PageSetup.Zoom = 49

$excel.CutCopyMode = 0
